# Add a "Sverige" (Sweden) total row beneath the per-region rows.
# File Reader for hospital and Intensive Care, added Critical stage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New summary row directly below the last region (row 21 -> row 22).
$ws.Range("A22").Value = "Sverige"
$ws.Range("B22").Value = 10300000

# Give the new number cell a thousands-separator format (#,##0).
$ws.Range("B22").NumberFormat = "#,##0"

# Move the active selection to the newly added cell.
$ws.Range("B22").Select()
